# Lab6 instructions paragraph edit:
#   " ... inherit from IdentityUser and your user needs to have some
#     properties of it's own."
# becomes four visually-identical but run-split pieces, "it's" is
# corrected to "its", and the (hidden) _GoBack bookmark - which used to
# sit at the very end of the paragraph - ends up sitting between "its"
# and " own." (i.e. where the last edit left off), matching the target
# revision's OOXML.

$d = $word.ActiveDocument
$apos = [char]0x2019

# --- locate "IdentityUser" -------------------------------------------------
$r1 = $d.Content
$r1.Find.Execute("IdentityUser", $false, $false, $false, $false, $false, `
                  $true, 1, $false, "", 0)
$iuStart = $r1.Start
$iuEnd   = $r1.End

# --- locate "it's" (the curly-quote apostrophe used in the source) --------
$r2 = $d.Content
$r2.Find.Execute("it" + $apos + "s own.", $false, $false, $false, $false, `
                  $false, $true, 1, $false, "", 0)
$itsStart = $r2.Start          # start of "it's"
$itsEnd   = $itsStart + 4       # "it's" is 4 characters long

# --- split "IdentityUser" into its own run ---------------------------------
# Toggling a character-format property on and back off leaves the visible
# formatting untouched but forces Word to break the enclosing run at the
# range boundaries, which is exactly the run-split the target shows.
$ruIU = $d.Range($iuStart, $iuEnd)
$ruIU.Bold = $true
$ruIU.Bold = $false

# --- split "it's" into its own run, then fix the apostrophe ---------------
$ruIts = $d.Range($itsStart, $itsEnd)
$ruIts.Bold = $true
$ruIts.Bold = $false

# Narrow text fix: drop the apostrophe ("it's" -> "its"). Only touch the
# "'s" tail so the already-isolated "IdentityUser" run is left alone.
$ruApos = $d.Range($itsStart + 2, $itsStart + 4)
$ruApos.Text = "s"

# The text edit above re-merges nearby same-format runs, so re-isolate
# the (now 3-character) "its" run as the very last touch to that span.
$itsNewEnd = $itsStart + 3
$ruItsFinal = $d.Range($itsStart, $itsNewEnd)
$ruItsFinal.Bold = $true
$ruItsFinal.Bold = $false

# --- relocate the hidden _GoBack bookmark ----------------------------------
# It used to sit at the end of the paragraph (after " own."); in the
# revised document it sits right after "its" (i.e. before " own."),
# matching where the edit actually happened.
$bmRange = $d.Range($itsNewEnd, $itsNewEnd)
$d.Bookmarks.Add("_GoBack", $bmRange)
